$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.21673309803009
$ws.Range("B1").Value = 1.577120661735535
$ws.Range("C1").Value = 2.111444711685181
$ws.Range("D1").Value = 6.12717866897583
$ws.Range("E1").Value = 3.026998519897461
